# Update computed results for each year's sheet with the latest server results.
$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("N2").Value = 7158.131594847751
$ws2025.Range("O2").Value = 6981.145263461227

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Value = 5707.815717280662
$ws2030.Range("I2").Value = 44492.05901988943
$ws2030.Range("L2").Value = 66334.06707325629
$ws2030.Range("M2").Value = 21991.42050229464
$ws2030.Range("N2").Value = 10598.18910437708
$ws2030.Range("O2").Value = 12082.01802941186

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 2927.360317916481
$ws2035.Range("B2").Value = 7940.887964949257
$ws2035.Range("E2").Value = 67179.99183625776
$ws2035.Range("I2").Value = 59530.75343380851
$ws2035.Range("L2").Value = 66334.06707325629
$ws2035.Range("M2").Value = 25547.11936466757
$ws2035.Range("N2").Value = 15125.63395528735
$ws2035.Range("O2").Value = 14763.87599448636

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Value = 2927.360317916481
$ws2040.Range("B2").Value = 7940.887964949257
$ws2040.Range("E2").Value = 67179.99183625776
$ws2040.Range("I2").Value = 59530.75343380851
$ws2040.Range("L2").Value = 66334.06707325629
$ws2040.Range("M2").Value = 25547.11936466757
$ws2040.Range("N2").Value = 15232.91015561297
$ws2040.Range("O2").Value = 14763.87599448636

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 6352.985609279765
$ws2045.Range("B2").Value = 7940.887964949257
$ws2045.Range("E2").Value = 67179.99183625776
$ws2045.Range("I2").Value = 59530.75343380851
$ws2045.Range("L2").Value = 66334.06707325629
$ws2045.Range("M2").Value = 25547.11936466757
$ws2045.Range("N2").Value = 15778.54231138309
$ws2045.Range("O2").Value = 17101.31291003395

$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("A2").Value = 6352.985609279765
$ws2050.Range("B2").Value = 7940.887964949257
$ws2050.Range("E2").Value = 67179.99183625776
$ws2050.Range("I2").Value = 59530.75343380851
$ws2050.Range("L2").Value = 66334.06707325629
$ws2050.Range("M2").Value = 25547.11936466757
$ws2050.Range("N2").Value = 15778.54231138309
$ws2050.Range("O2").Value = 17101.31291003395
